$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = "OK SCRIPT"
$ws.Range("I2").Font.Bold = $true
$ws.Range("I2").HorizontalAlignment = -4108
$ws.Range("I2").VerticalAlignment = -4108

$ws.Range("F3").Copy()
$ws.Range("I3:I30").PasteSpecial(-4122)

$rng = $ws.Range("I3:I29")
$rng.Value = "X"

Write-Host "done"
